$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Hovedbruker" (main user) row - change the password shown in columns
# C, F, I and L of row 4 from the old shared password ("Maine1953",
# shared-string 22) to the new one ("Roma1995", shared-string 24) so it
# matches the "Eier" (owner) row above it.
$ws.Range("C4").Value = "Roma1995"
$ws.Range("F4").Value = "Roma1995"
$ws.Range("I4").Value = "Roma1995"
$ws.Range("L4").Value = "Roma1995"

# Move/extend the on-screen selection to K15:K16 (as last left by the author).
$ws.Range("K15:K16").Select()
